# Version 0.3 - 23.09.2024 - Change of substructure and code
#
# - Adds four new sheets (product_origin, product_fractions, product_quality,
#   product_amount) after the existing product_data sheet.
# - Adds a new response row (row 9) to product_data.
# - Clears the stray empty string cells left over in company_data row 62.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. company_data (sheet2): row 62 had a handful of cells holding empty
#    strings left behind by the export script; blank them out completely.
# ---------------------------------------------------------------------------
$company = $wb.Worksheets.Item("company_data")
foreach ($col in @("C", "D", "E", "F", "G", "H", "K")) {
    $company.Range($col + "62").Value = ""
}

# ---------------------------------------------------------------------------
# 2. product_data (sheet3): append a new submission as row 9.
# ---------------------------------------------------------------------------
$product = $wb.Worksheets.Item("product_data")
$product.Range("A9").Value = 1
$product.Range("B9").Value = "2024-09-16 09:30:27"
$product.Range("C9").Value = "['Kunststoff – ABS', 'Kunststoff – ASA']"
$product.Range("D9").Value = "['ABS-Test', 'ASA-Test']"
$product.Range("E9").Value = "[90.0, 10.0]"

# ---------------------------------------------------------------------------
# Helper: add a brand-new worksheet at the end of the workbook.
# ---------------------------------------------------------------------------
function Add-SheetAtEnd([string]$name) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
    $newSheet.Name = $name
    return $newSheet
}

# ---------------------------------------------------------------------------
# 3. product_origin (new sheet4)
# ---------------------------------------------------------------------------
$origin = Add-SheetAtEnd "product_origin"

$originRows = @(
    @{ A=1; B="2024-09-16 12:10:56"; C="Post-Consumer (PC) – getrennte Sammlung"; D="Test"; E="lokal als Bringsystem"; F="12 15 15*" },
    @{ A=1; B="2024-09-16 12:12:05"; C="Post-Consumer (PC) – getrennte Sammlung"; D="Test"; E="regional" },
    @{ A=1; B="2024-09-16 12:15:18"; C="Post-Industrial (PI)"; D="Test" },
    @{ A=1; B="2024-09-16 12:19:11"; C="Post-Industrial (PI)"; D="Test" },
    @{ A=1; B="2024-09-16 12:19:24"; C="Post-Consumer (PC) – getrennte Sammlung"; D="Test"; E="lokal als Bringsystem"; F="15215" },
    @{ A=1; B="2024-09-16 13:17:48"; C="Post-Industrial (PI)"; D="Test" },
    @{ A=1; B="2024-09-17 10:57:36"; C="Post-Industrial (PI)" }
)

# F5 ("15215") looks numeric but must stay text, matching the source export.
$origin.Range("F5").NumberFormat = "@"

$r = 1
foreach ($row in $originRows) {
    foreach ($col in $row.Keys) {
        $origin.Range($col + $r).Value = $row[$col]
    }
    $r++
}

# ---------------------------------------------------------------------------
# 4. product_fractions (new sheet5)
# ---------------------------------------------------------------------------
$fractions = Add-SheetAtEnd "product_fractions"

$fractionRows = @(
    @{ A=1; B="2024-09-16 12:11:21"; C="['Kunststoff – ABS', 'Kunststoff – ASA']"; D="['', '']"; E="[0.0, 0.0]" },
    @{ A=1; B="2024-09-16 13:17:45"; C="['Kunststoff – ABS', 'Kunststoff – ABS']"; D="['', '']"; E="[0.0, 0.0]" },
    @{ A=1; B="2024-09-17 10:57:24"; C="['Kunststoff – ABS', 'Kunststoff – ASA']"; D="['', '']"; E="[0.0, 0.0]" }
)

$r = 1
foreach ($row in $fractionRows) {
    foreach ($col in $row.Keys) {
        $fractions.Range($col + $r).Value = $row[$col]
    }
    $r++
}

# ---------------------------------------------------------------------------
# 5. product_quality (new sheet6)
# ---------------------------------------------------------------------------
$quality = Add-SheetAtEnd "product_quality"

$qualityRows = @(
    @{ A=1; B="2024-09-16 14:03:27"; C="Ja"; D="blau"; E=0; F="keine"; H="['Antibeschlagmittel', 'Antibeschlagmittel']"; I="[0.0, 0.0]" },
    @{ A=1; B="2024-09-16 14:03:50"; C="Ja"; D="blau"; E=0; F="keine"; H="['Antibeschlagmittel', 'Antibeschlagmittel']"; I="[0.0, 0.0]" },
    @{ A=1; B="2024-09-16 14:08:01"; C="Ja"; D="blau"; E=0; F="keine"; H="['Antibeschlagmittel', 'Antibeschlagmittel', 'Antibeschlagmittel']"; I="[0.0, 0.0, 0.0]" },
    @{ A=1; B="2024-09-17 11:19:33"; E=100; H="[['Antioxidant'], ['Gleitmittel']]"; I="[]" },
    @{ A=1; B="2024-09-17 11:19:47"; E=100; H="[['Antioxidant', 'Biozide', 'Gleitmittel'], ['Gleitmittel', 'Antioxidant']]"; I="[]" },
    @{ A=1; B="2024-09-19 11:22:09"; E=100; H="[[], [], []]"; I="[['Glimmer'], ['Ruß'], ['Glaskugeln', 'Glimmer', 'Glasfasern']]" },
    @{ A=1; B="2024-09-19 11:26:11"; E=100; H="[['Gleitmittel'], [], []]"; I="[[], [], []]" },
    @{ A=1; B="2024-09-19 11:28:15"; C=""; D=""; E=100; F=""; G=""; H="[[], [], []]"; I="[[], [], []]" }
)

$r = 1
foreach ($row in $qualityRows) {
    foreach ($col in $row.Keys) {
        $quality.Range($col + $r).Value = $row[$col]
    }
    $r++
}

# ---------------------------------------------------------------------------
# 6. product_amount (new sheet7)
# ---------------------------------------------------------------------------
$amount = Add-SheetAtEnd "product_amount"

$amountRows = @(
    @{ A=1; B="2024-09-16 14:39:13"; C=0; D=0; E="Tag" },
    @{ A=1; B="2024-09-16 14:40:18"; C=20; D=20; E="Tag"; F="Test`n" },
    @{ A=1; B="2024-09-16 14:40:36"; C=0; D=0; E="Tag" },
    @{ A=1; B="2024-09-16 14:41:38"; C=20; D=10; E="Woche"; F="12" }
)

# F4 ("12") looks numeric but must stay text, matching the source export.
$amount.Range("F4").NumberFormat = "@"

$r = 1
foreach ($row in $amountRows) {
    foreach ($col in $row.Keys) {
        $amount.Range($col + $r).Value = $row[$col]
    }
    $r++
}

# ---------------------------------------------------------------------------
# Leave the view pointed at the first sheet, as in the original workbook.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("contact_data").Activate()
